$d = $word.ActiveDocument

# The document's first paragraph originally holds two runs:
#   "**ID__AFFARS_mp_5342_902_topic_1__ID**"  +  " " (trailing space run)
# The edit collapses this into a single run with updated placeholder text,
# and also adds paragraph border "space" (no visible lines) plus a new
# left indent.

$p1 = $d.Paragraphs(1)

# --- Update paragraph formatting -----------------------------------------
# w:ind w:left="120" -> w:ind w:left="225"  (LeftIndent is expressed in
# points; 1 twip = 1/20 pt, so 225 twips = 11.25 pt)
$p1.Range.ParagraphFormat.LeftIndent = 225 / 20.0

# Add <w:pBdr><w:top w:space="5"/><w:left w:space="5"/><w:bottom w:space="5"/>
#            <w:right w:space="5"/></w:pBdr>
# (space-only paragraph border, no visible line)
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# --- Update the paragraph's text content ---------------------------------
# Remove the trailing standalone space run, then rewrite the remaining
# run's text in place so only a single run remains in the paragraph.
$full = $p1.Range
$full.End = $full.End - 1               # exclude paragraph mark
$full.Text = "**ID__AFFARS_MP5342_902__ID**"
